$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Efnb2"
$ws.Cells.Item(2, 3).Value = "Epha3"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 45.76217133333333
$ws.Cells.Item(2, 8).Value = 137.286514
$ws.Cells.Item(2, 9).Value = 0.6763939203605134
$ws.Cells.Item(2, 10).Value = 0.6763939203605135
$ws.Cells.Item(2, 11).Value = 2
$ws.Cells.Item(2, 12).Value = 0.6666666666666666
$ws.Cells.Item(2, 13).Value = 0.07579599999999999
$ws.Cells.Item(2, 14).Value = 0.227388
$ws.Cells.Item(2, 15).Value = 0.001780200955210419
$ws.Cells.Item(2, 16).Value = 0.001780200955210419
$ws.Cells.Item(2, 17).Value = 3.468589538381333
$ws.Cells.Item(2, 18).Value = 31.217305845432
$ws.Cells.Item(2, 19).Value = 0.001204117103124306
$ws.Cells.Item(2, 20).Value = 0.001204117103124306

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Efnb2"
$ws.Cells.Item(3, 3).Value = "Epha3"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 45.76217133333333
$ws.Cells.Item(3, 8).Value = 137.286514
$ws.Cells.Item(3, 9).Value = 0.6763939203605134
$ws.Cells.Item(3, 10).Value = 0.6763939203605135
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 40.12734033333333
$ws.Cells.Item(3, 14).Value = 120.382021
$ws.Cells.Item(3, 15).Value = 0.9424604146848589
$ws.Cells.Item(3, 16).Value = 0.9424604146848587
$ws.Cells.Item(3, 17).Value = 1836.314223484977
$ws.Cells.Item(3, 18).Value = 16526.8280113648
$ws.Cells.Item(3, 19).Value = 0.6374744946732869
$ws.Cells.Item(3, 20).Value = 0.6374744946732869

# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Efnb2"
$ws.Cells.Item(4, 3).Value = "Epha3"
$ws.Cells.Item(4, 4).Value = "MuSCs"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 45.76217133333333
$ws.Cells.Item(4, 8).Value = 137.286514
$ws.Cells.Item(4, 9).Value = 0.6763939203605134
$ws.Cells.Item(4, 10).Value = 0.6763939203605135
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 2.344072666666666
$ws.Cells.Item(4, 14).Value = 7.032217999999999
$ws.Cells.Item(4, 15).Value = 0.05505462557763778
$ws.Cells.Item(4, 16).Value = 0.05505462557763778
$ws.Cells.Item(4, 17).Value = 107.2698549897836
$ws.Cells.Item(4, 18).Value = 965.428694908052
$ws.Cells.Item(4, 19).Value = 0.03723861402843861
$ws.Cells.Item(4, 20).Value = 0.03723861402843862

# Row 5
$ws.Cells.Item(5, 1).Value = "ECs"
$ws.Cells.Item(5, 2).Value = "Efnb2"
$ws.Cells.Item(5, 3).Value = "Epha3"
$ws.Cells.Item(5, 4).Value = "Resolving-Mac"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 45.76217133333333
$ws.Cells.Item(5, 8).Value = 137.286514
$ws.Cells.Item(5, 9).Value = 0.6763939203605134
$ws.Cells.Item(5, 10).Value = 0.6763939203605135
$ws.Cells.Item(5, 11).Value = 2
$ws.Cells.Item(5, 12).Value = 0.6666666666666666
$ws.Cells.Item(5, 13).Value = 0.03000666666666667
$ws.Cells.Item(5, 14).Value = 0.09002
$ws.Cells.Item(5, 15).Value = 0.0007047587822930054
$ws.Cells.Item(5, 16).Value = 0.0007047587822930053
$ws.Cells.Item(5, 17).Value = 1.373170221142222
$ws.Cells.Item(5, 18).Value = 12.35853199028
$ws.Cells.Item(5, 19).Value = 0.0004766945556636675
$ws.Cells.Item(5, 20).Value = 0.0004766945556636675

# Row 6
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Efnb2"
$ws.Cells.Item(6, 3).Value = "Epha3"
$ws.Cells.Item(6, 4).Value = "ECs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 9.395935333333332
$ws.Cells.Item(6, 8).Value = 28.187806
$ws.Cells.Item(6, 9).Value = 0.1388778842960613
$ws.Cells.Item(6, 10).Value = 0.1388778842960613
$ws.Cells.Item(6, 11).Value = 2
$ws.Cells.Item(6, 12).Value = 0.6666666666666666
$ws.Cells.Item(6, 13).Value = 0.07579599999999999
$ws.Cells.Item(6, 14).Value = 0.227388
$ws.Cells.Item(6, 15).Value = 0.001780200955210419
$ws.Cells.Item(6, 16).Value = 0.001780200955210419
$ws.Cells.Item(6, 17).Value = 0.7121743145253331
$ws.Cells.Item(6, 18).Value = 6.409568830727999
$ws.Cells.Item(6, 19).Value = 0.0002472305422814503
$ws.Cells.Item(6, 20).Value = 0.0002472305422814504

# Row 7
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Efnb2"
$ws.Cells.Item(7, 3).Value = "Epha3"
$ws.Cells.Item(7, 4).Value = "FAPs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 9.395935333333332
$ws.Cells.Item(7, 8).Value = 28.187806
$ws.Cells.Item(7, 9).Value = 0.1388778842960613
$ws.Cells.Item(7, 10).Value = 0.1388778842960613
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 40.12734033333333
$ws.Cells.Item(7, 14).Value = 120.382021
$ws.Cells.Item(7, 15).Value = 0.9424604146848589
$ws.Cells.Item(7, 16).Value = 0.9424604146848587
$ws.Cells.Item(7, 17).Value = 377.0338948706583
$ws.Cells.Item(7, 18).Value = 3393.305053835926
$ws.Cells.Item(7, 19).Value = 0.1308869084242218
$ws.Cells.Item(7, 20).Value = 0.1308869084242218

# Row 8
$ws.Cells.Item(8, 1).Value = "FAPs"
$ws.Cells.Item(8, 2).Value = "Efnb2"
$ws.Cells.Item(8, 3).Value = "Epha3"
$ws.Cells.Item(8, 4).Value = "MuSCs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 9.395935333333332
$ws.Cells.Item(8, 8).Value = 28.187806
$ws.Cells.Item(8, 9).Value = 0.1388778842960613
$ws.Cells.Item(8, 10).Value = 0.1388778842960613
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 2.344072666666666
$ws.Cells.Item(8, 14).Value = 7.032217999999999
$ws.Cells.Item(8, 15).Value = 0.05505462557763778
$ws.Cells.Item(8, 16).Value = 0.05505462557763778
$ws.Cells.Item(8, 17).Value = 22.02475519263422
$ws.Cells.Item(8, 18).Value = 198.222796733708
$ws.Cells.Item(8, 19).Value = 0.007645869920934157
$ws.Cells.Item(8, 20).Value = 0.007645869920934158

# Row 9
$ws.Cells.Item(9, 1).Value = "FAPs"
$ws.Cells.Item(9, 2).Value = "Efnb2"
$ws.Cells.Item(9, 3).Value = "Epha3"
$ws.Cells.Item(9, 4).Value = "Resolving-Mac"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 9.395935333333332
$ws.Cells.Item(9, 8).Value = 28.187806
$ws.Cells.Item(9, 9).Value = 0.1388778842960613
$ws.Cells.Item(9, 10).Value = 0.1388778842960613
$ws.Cells.Item(9, 11).Value = 2
$ws.Cells.Item(9, 12).Value = 0.6666666666666666
$ws.Cells.Item(9, 13).Value = 0.03000666666666667
$ws.Cells.Item(9, 14).Value = 0.09002
$ws.Cells.Item(9, 15).Value = 0.0007047587822930054
$ws.Cells.Item(9, 16).Value = 0.0007047587822930053
$ws.Cells.Item(9, 17).Value = 0.2819406995688888
$ws.Cells.Item(9, 18).Value = 2.53746629612
$ws.Cells.Item(9, 19).Value = 0.00009787540862392106
$ws.Cells.Item(9, 20).Value = 0.00009787540862392106

# Row 10
$ws.Cells.Item(10, 1).Value = "MuSCs"
$ws.Cells.Item(10, 2).Value = "Efnb2"
$ws.Cells.Item(10, 3).Value = "Epha3"
$ws.Cells.Item(10, 4).Value = "ECs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 12.29750866666667
$ws.Cells.Item(10, 8).Value = 36.892526
$ws.Cells.Item(10, 9).Value = 0.1817649787009828
$ws.Cells.Item(10, 10).Value = 0.1817649787009828
$ws.Cells.Item(10, 11).Value = 2
$ws.Cells.Item(10, 12).Value = 0.6666666666666666
$ws.Cells.Item(10, 13).Value = 0.07579599999999999
$ws.Cells.Item(10, 14).Value = 0.227388
$ws.Cells.Item(10, 15).Value = 0.001780200955210419
$ws.Cells.Item(10, 16).Value = 0.001780200955210419
$ws.Cells.Item(10, 17).Value = 0.9321019668986665
$ws.Cells.Item(10, 18).Value = 8.388917702088
$ws.Cells.Item(10, 19).Value = 0.0003235781887072909
$ws.Cells.Item(10, 20).Value = 0.0003235781887072909

# Row 11
$ws.Cells.Item(11, 1).Value = "MuSCs"
$ws.Cells.Item(11, 2).Value = "Efnb2"
$ws.Cells.Item(11, 3).Value = "Epha3"
$ws.Cells.Item(11, 4).Value = "FAPs"
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 12.29750866666667
$ws.Cells.Item(11, 8).Value = 36.892526
$ws.Cells.Item(11, 9).Value = 0.1817649787009828
$ws.Cells.Item(11, 10).Value = 0.1817649787009828
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 13).Value = 40.12734033333333
$ws.Cells.Item(11, 14).Value = 120.382021
$ws.Cells.Item(11, 15).Value = 0.9424604146848589
$ws.Cells.Item(11, 16).Value = 0.9424604146848587
$ws.Cells.Item(11, 17).Value = 493.4663155194496
$ws.Cells.Item(11, 18).Value = 4441.196839675046
$ws.Cells.Item(11, 19).Value = 0.1713062972017128
$ws.Cells.Item(11, 20).Value = 0.1713062972017127

# Row 12
$ws.Cells.Item(12, 1).Value = "MuSCs"
$ws.Cells.Item(12, 2).Value = "Efnb2"
$ws.Cells.Item(12, 3).Value = "Epha3"
$ws.Cells.Item(12, 4).Value = "MuSCs"
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 12.29750866666667
$ws.Cells.Item(12, 8).Value = 36.892526
$ws.Cells.Item(12, 9).Value = 0.1817649787009828
$ws.Cells.Item(12, 10).Value = 0.1817649787009828
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 12).Value = 1
$ws.Cells.Item(12, 13).Value = 2.344072666666666
$ws.Cells.Item(12, 14).Value = 7.032217999999999
$ws.Cells.Item(12, 15).Value = 0.05505462557763778
$ws.Cells.Item(12, 16).Value = 0.05505462557763778
$ws.Cells.Item(12, 17).Value = 28.82625393362978
$ws.Cells.Item(12, 18).Value = 259.436285402668
$ws.Cells.Item(12, 19).Value = 0.01000700284550991
$ws.Cells.Item(12, 20).Value = 0.01000700284550991

# Row 13
$ws.Cells.Item(13, 1).Value = "MuSCs"
$ws.Cells.Item(13, 2).Value = "Efnb2"
$ws.Cells.Item(13, 3).Value = "Epha3"
$ws.Cells.Item(13, 4).Value = "Resolving-Mac"
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 12.29750866666667
$ws.Cells.Item(13, 8).Value = 36.892526
$ws.Cells.Item(13, 9).Value = 0.1817649787009828
$ws.Cells.Item(13, 10).Value = 0.1817649787009828
$ws.Cells.Item(13, 11).Value = 2
$ws.Cells.Item(13, 12).Value = 0.6666666666666666
$ws.Cells.Item(13, 13).Value = 0.03000666666666667
$ws.Cells.Item(13, 14).Value = 0.09002
$ws.Cells.Item(13, 15).Value = 0.0007047587822930054
$ws.Cells.Item(13, 16).Value = 0.0007047587822930053
$ws.Cells.Item(13, 17).Value = 0.3690072433911111
$ws.Cells.Item(13, 18).Value = 3.321065190520001
$ws.Cells.Item(13, 19).Value = 0.0001281004650528187
$ws.Cells.Item(13, 20).Value = 0.0001281004650528187

# Row 14
$ws.Cells.Item(14, 1).Value = "Resolving-Mac"
$ws.Cells.Item(14, 2).Value = "Efnb2"
$ws.Cells.Item(14, 3).Value = "Epha3"
$ws.Cells.Item(14, 4).Value = "ECs"
$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 6).Value = 1
$ws.Cells.Item(14, 7).Value = 0.2004796666666666
$ws.Cells.Item(14, 8).Value = 0.601439
$ws.Cells.Item(14, 9).Value = 0.002963216642442438
$ws.Cells.Item(14, 10).Value = 0.002963216642442439
$ws.Cells.Item(14, 11).Value = 2
$ws.Cells.Item(14, 12).Value = 0.6666666666666666
$ws.Cells.Item(14, 13).Value = 0.07579599999999999
$ws.Cells.Item(14, 14).Value = 0.227388
$ws.Cells.Item(14, 15).Value = 0.001780200955210419
$ws.Cells.Item(14, 16).Value = 0.001780200955210419
$ws.Cells.Item(14, 17).Value = 0.01519555681466666
$ws.Cells.Item(14, 18).Value = 0.136760011332
$ws.Cells.Item(14, 19).Value = 0.000005275121097371438
$ws.Cells.Item(14, 20).Value = 0.000005275121097371439

# Row 15
$ws.Cells.Item(15, 1).Value = "Resolving-Mac"
$ws.Cells.Item(15, 2).Value = "Efnb2"
$ws.Cells.Item(15, 3).Value = "Epha3"
$ws.Cells.Item(15, 4).Value = "FAPs"
$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 6).Value = 1
$ws.Cells.Item(15, 7).Value = 0.2004796666666666
$ws.Cells.Item(15, 8).Value = 0.601439
$ws.Cells.Item(15, 9).Value = 0.002963216642442438
$ws.Cells.Item(15, 10).Value = 0.002963216642442439
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 12).Value = 1
$ws.Cells.Item(15, 13).Value = 40.12734033333333
$ws.Cells.Item(15, 14).Value = 120.382021
$ws.Cells.Item(15, 15).Value = 0.9424604146848589
$ws.Cells.Item(15, 16).Value = 0.9424604146848587
$ws.Cells.Item(15, 17).Value = 8.044715814246553
$ws.Cells.Item(15, 18).Value = 72.402442328219
$ws.Cells.Item(15, 19).Value = 0.002792714385637375
$ws.Cells.Item(15, 20).Value = 0.002792714385637375

# Row 16
$ws.Cells.Item(16, 1).Value = "Resolving-Mac"
$ws.Cells.Item(16, 2).Value = "Efnb2"
$ws.Cells.Item(16, 3).Value = "Epha3"
$ws.Cells.Item(16, 4).Value = "MuSCs"
$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 6).Value = 1
$ws.Cells.Item(16, 7).Value = 0.2004796666666666
$ws.Cells.Item(16, 8).Value = 0.601439
$ws.Cells.Item(16, 9).Value = 0.002963216642442438
$ws.Cells.Item(16, 10).Value = 0.002963216642442439
$ws.Cells.Item(16, 11).Value = 3
$ws.Cells.Item(16, 12).Value = 1
$ws.Cells.Item(16, 13).Value = 2.344072666666666
$ws.Cells.Item(16, 14).Value = 7.032217999999999
$ws.Cells.Item(16, 15).Value = 0.05505462557763778
$ws.Cells.Item(16, 16).Value = 0.05505462557763778
$ws.Cells.Item(16, 17).Value = 0.4699389068557777
$ws.Cells.Item(16, 18).Value = 4.229450161701999
$ws.Cells.Item(16, 19).Value = 0.0001631387827550934
$ws.Cells.Item(16, 20).Value = 0.0001631387827550934

# Row 17
$ws.Cells.Item(17, 1).Value = "Resolving-Mac"
$ws.Cells.Item(17, 2).Value = "Efnb2"
$ws.Cells.Item(17, 3).Value = "Epha3"
$ws.Cells.Item(17, 4).Value = "Resolving-Mac"
$ws.Cells.Item(17, 5).Value = 3
$ws.Cells.Item(17, 6).Value = 1
$ws.Cells.Item(17, 7).Value = 0.2004796666666666
$ws.Cells.Item(17, 8).Value = 0.601439
$ws.Cells.Item(17, 9).Value = 0.002963216642442438
$ws.Cells.Item(17, 10).Value = 0.002963216642442439
$ws.Cells.Item(17, 11).Value = 2
$ws.Cells.Item(17, 12).Value = 0.6666666666666666
$ws.Cells.Item(17, 13).Value = 0.03000666666666667
$ws.Cells.Item(17, 14).Value = 0.09002
$ws.Cells.Item(17, 15).Value = 0.0007047587822930054
$ws.Cells.Item(17, 16).Value = 0.0007047587822930053
$ws.Cells.Item(17, 17).Value = 0.006015726531111111
$ws.Cells.Item(17, 18).Value = 0.05414153878
$ws.Cells.Item(17, 19).Value = 0.000002088352952598101
$ws.Cells.Item(17, 20).Value = 0.000002088352952598101
